$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 2.250477639167457
$ws.Range("C16").Value = 44.89971791433486
$ws.Range("D16").Value = 17.74952915753526
$ws.Range("E16").Value = 1244.949559337892
